$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (H1) so the new headers
# match the bold/centered/bordered header formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-6 for columns I (I0) and J (IF)
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 9
